$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.263.96"
$ws.Range("E2").Value = "  +1.58%  "
$ws.Range("D3").Value = "3.893.24"
$ws.Range("E3").Value = "  -0.39%  "
$ws.Range("D4").Value = "'0.997"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.48%  "
$ws.Range("D5").Value = "'522.87"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +8.15%  "
$ws.Range("D6").Value = "'143.18"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.32%  "
$ws.Range("D7").Value = "'0.610"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.52%  "
$ws.Range("D8").Value = "'0.998"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  -2.53%  "
$ws.Range("E10").Value = "  +2.50%  "
$ws.Range("D11").Value = "'0.0000333"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.49%  "
$ws.Range("D12").Value = "'41.89"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.73%  "
$ws.Range("D13").Value = "4.511.29"
$ws.Range("E13").Value = "  -0.65%  "
$ws.Range("D14").Value = "'10.21"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.88%  "
$ws.Range("D15").Value = "3.890.99"
$ws.Range("E15").Value = "  -0.53%  "
$ws.Range("B16").Value = "Uniswap"
$ws.Range("C16").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D16").Value = "'13.96"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.65%  "
$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D17").Value = "'0.135"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.87%  "
$ws.Range("D18").Value = "'1.21"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +6.66%  "
$ws.Range("E19").Value = "  -2.92%  "
$ws.Range("D20").Value = "69.142.17"
$ws.Range("E20").Value = "  +1.20%  "
$ws.Range("D21").Value = "'425.22"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.10%  "
$ws.Range("E22").Value = "  -5.68%  "
$ws.Range("B23").Value = "Litecoin"
$ws.Range("C23").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D23").Value = "'88.17"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.23%  "
$ws.Range("B24").Value = "InternetComputer(DFINITY)"
$ws.Range("C24").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D24").Value = "'14.09"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -6.13%  "
$ws.Range("D25").Value = "'3.98"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +7.18%  "
$ws.Range("D26").Value = "'11.34"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.59%  "
$ws.Range("D27").Value = "'10.48"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.45%  "
$ws.Range("D28").Value = "'36.14"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.78%  "
$ws.Range("D29").Value = "'693.51"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.23%  "
$ws.Range("D30").Value = "'13.09"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.58%  "
$ws.Range("E31").Value = "  -4.60%  "
$ws.Range("D32").Value = "'2.81"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.05%  "
$ws.Range("D33").Value = "'67.81"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +11.81%  "
$ws.Range("D34").Value = "'0.433"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +8.00%  "
$ws.Range("D35").Value = "'5.87"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.69%  "
$ws.Range("D36").Value = "'39.73"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.06%  "
$ws.Range("E37").Value = "  -4.78%  "
$ws.Range("E38").Value = "  +2.43%  "
$ws.Range("D39").Value = "'0.997"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.16%  "
$ws.Range("E40").Value = "  -0.23%  "
$ws.Range("E41").Value = "  -3.00%  "
$ws.Range("D42").Value = "'3.08"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.54%  "
$ws.Range("D43").Value = "'2.77"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -8.80%  "
$ws.Range("D44").Value = "'3.00"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.85%  "
$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D45").Value = "'0.140"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.92%  "
$ws.Range("B46").Value = "ApeXProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D46").Value = "'3.32"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.82%  "
$ws.Range("D47").Value = "'3.00"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.83%  "
$ws.Range("D48").Value = "'26.93"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +6.68%  "
$ws.Range("B49").Value = "Maker"
$ws.Range("C49").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D49").Value = "2.727.78"
$ws.Range("E49").Value = "  +11.56%  "
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "0.0₆0340"
$ws.Range("E50").Value = "  +2.25%  "
$ws.Range("E51").Value = "  -4.60%  "
